$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" '61.542.58'
Set-TextValue $ws "E2" '  +0.93%  '

Set-TextValue $ws "D3" '3.387.56'
Set-TextValue $ws "E3" '  +0.49%  '

Set-TextValue $ws "E4" '  -0.01%  '

Set-TextValue $ws "D5" '576.72'
Set-TextValue $ws "E5" '  +1.04%  '

Set-TextValue $ws "D6" '137.09'
Set-TextValue $ws "E6" '  +1.19%  '

Set-TextValue $ws "E7" '  -0.01%  '

Set-TextValue $ws "D8" '3.388.20'
Set-TextValue $ws "E8" '  +0.57%  '

Set-TextValue $ws "D9" '0.473'
Set-TextValue $ws "E9" '  -0.58%  '

Set-TextValue $ws "D10" '7.50'
Set-TextValue $ws "E10" '  -1.00%  '

Set-TextValue $ws "D11" '0.125'
Set-TextValue $ws "E11" '  +2.37%  '

Set-TextValue $ws "E12" '  +0.32%  '

Set-TextValue $ws "D13" '3.962.48'
Set-TextValue $ws "E13" '  +0.28%  '

Set-TextValue $ws "E14" '  +1.54%  '

Set-TextValue $ws "E15" '  +1.67%  '

Set-TextValue $ws "D16" '3.388.84'
Set-TextValue $ws "E16" '  +0.47%  '

Set-TextValue $ws "D17" '25.65'
Set-TextValue $ws "E17" '  +2.45%  '

Set-TextValue $ws "D18" '61.649.61'
Set-TextValue $ws "E18" '  +0.82%  '

Set-TextValue $ws "D19" '14.16'
Set-TextValue $ws "E19" '  +1.07%  '

Set-TextValue $ws "D20" '9.47'
Set-TextValue $ws "E20" '  +0.50%  '

Set-TextValue $ws "D21" '5.80'
Set-TextValue $ws "E21" '  +0.58%  '

Set-TextValue $ws "D22" '377.64'
Set-TextValue $ws "E22" '  +1.07%  '

Set-TextValue $ws "D23" '0.559'
Set-TextValue $ws "E23" '  -1.21%  '

Set-TextValue $ws "D24" '3.524.00'
Set-TextValue $ws "E24" '  +0.54%  '

Set-TextValue $ws "E25" '  -0.09%  '

Set-TextValue $ws "B26" 'Litecoin'
Set-TextValue $ws "C26" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws "D26" '71.04'
Set-TextValue $ws "E26" '  +0.75%  '

Set-TextValue $ws "B27" 'PEPE'
Set-TextValue $ws "C27" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws "D27" '0.0000124'
Set-TextValue $ws "E27" '  +6.47%  '

Set-TextValue $ws "D28" '1.74'
Set-TextValue $ws "E28" '  +5.22%  '

Set-TextValue $ws "D29" '7.58'
Set-TextValue $ws "E29" '  -1.31%  '

Set-TextValue $ws "D30" '0.997'
Set-TextValue $ws "E30" '  -0.30%  '

Set-TextValue $ws "B31" 'InternetComputer(DFINITY)'
Set-TextValue $ws "C31" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws "D31" '8.18'
Set-TextValue $ws "E31" '  +1.05%  '

Set-TextValue $ws "B32" 'Kaspa'
Set-TextValue $ws "C32" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D32" '0.159'
Set-TextValue $ws "E32" '  +3.41%  '

Set-TextValue $ws "E33" '  +0.88%  '

Set-TextValue $ws "E34" '  +0.05%  '

Set-TextValue $ws "D35" '23.35'
Set-TextValue $ws "E35" '  +0.17%  '

Set-TextValue $ws "E36" '  -3.47%  '

Set-TextValue $ws "E37" '  -0.05%  '

Set-TextValue $ws "D38" '6.84'
Set-TextValue $ws "E38" '  -1.04%  '

Set-TextValue $ws "D39" '164.48'
Set-TextValue $ws "E39" '  +0.85%  '

Set-TextValue $ws "D40" '0.0783'
Set-TextValue $ws "E40" '  -0.47%  '

Set-TextValue $ws "B41" 'ONDO'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws "D41" '1.24'
Set-TextValue $ws "E41" '  +2.98%  '

Set-TextValue $ws "B42" 'FirstDigitalUSD'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws "D42" '1.00'
Set-TextValue $ws "E42" '  +0.07%  '

Set-TextValue $ws "B43" 'Mantle'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws "D43" '0.780'
Set-TextValue $ws "E43" '  +2.87%  '

Set-TextValue $ws "E44" '  +7.47%  '

Set-TextValue $ws "B45" 'Filecoin'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D45" '4.41'
Set-TextValue $ws "E45" '  +0.00%  '

Set-TextValue $ws "B46" 'EnergySwap'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D46" '24.87'
Set-TextValue $ws "E46" '  +7.82%  '

Set-TextValue $ws "D47" '41.29'
Set-TextValue $ws "E47" '  +0.23%  '

Set-TextValue $ws "D48" '6.85'
Set-TextValue $ws "E48" '  -1.86%  '

Set-TextValue $ws "D49" '22.76'
Set-TextValue $ws "E49" '  -0.78%  '

Set-TextValue $ws "D50" '2.331.33'
Set-TextValue $ws "E50" '  +5.37%  '

Set-TextValue $ws "E51" '  +1.86%  '
